$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")

$ws.Range("A7").Value = "test_search_result"
$ws.Range("B7").Value = "XPath"
$ws.Range("C7").Value = "//a[@title='searchText']"

$ws.Range("C8").Value = "//span[@class='-link--channel-name']"
$ws.Range("A8").Value = "nav_bar_stackoverflow_menu"
$ws.Range("B8").Value = "XPath"

$ws.Range("C9").Value = "//a[@class='pl8 js-gps-track nav-links--link']"
$ws.Range("A9").Value = "nav_bar_home_menu"
$ws.Range("B9").Value = "XPath"

$ws.Range("A9").Select() | Out-Null
